$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in E2:E9 from 50 to 70
$ws.Range("E2:E9").Value = 70

# Update the selection to E2:E9 with active cell E2
$ws.Range("E2:E9").Select()
